$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") for rows 2-29: update date serial 45402 -> 45406
for ($row = 2; $row -le 29; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45402) {
        $cell.Value2 = 45406
    }
}
